# Cotações atualizadas - 2025-10-09
# Append a new row (35) with the quotes for 2025-10-09 (serial date 45939),
# mirroring the existing rows' layout: date in column A, values as text
# strings (Portuguese comma-decimal format) in columns B:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 35
$prevRow = $newRow - 1

# Date value for 2025-10-09 (Excel serial 45939)
$ws.Range("A$newRow").Value = 45939
# Match the date formatting used by the previous row (style s="2")
$ws.Range("A$newRow").NumberFormat = $ws.Range("A$prevRow").NumberFormat

$ws.Range("B$newRow").Value = "21,7372"
$ws.Range("C$newRow").Value = "15,6392"
$ws.Range("D$newRow").Value = "15,5076"
$ws.Range("E$newRow").Value = "15,5076"
